# Insert a new data row at row 12 (pushing the existing rows 12-72 down to 13-73)
# and populate it with the new weekly price-report entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value2 = 10
$ws.Range("B12").Value2 = "Vega Modelo de Temuco"
$ws.Range("C12").Value2 = "La Araucanía"
$ws.Range("D12").Value2 = 44537
$ws.Range("E12").Value2 = 9
$ws.Range("F12").Value2 = 100112022
$ws.Range("G12").Value2 = "Arveja Verde"
$ws.Range("H12").Value2 = "Sin especificar"
$ws.Range("I12").Value2 = "Primera"
$ws.Range("J12").Value2 = 95
$ws.Range("K12").Value2 = 13000
$ws.Range("L12").Value2 = 13000
$ws.Range("M12").Value2 = 13000
$ws.Range("N12").Value2 = "`$/saco 25 kilos"
$ws.Range("O12").Value2 = "Región de La Araucanía"
$ws.Range("P12").Value2 = 520
$ws.Range("Q12").Value2 = 25
$ws.Range("R12").Value2 = "Hortaliza"
